$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("C14").Value = 410
$ws.Range("J31").Select()
